$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 339
$newDate = 45186

# Column C ("Förändrad") gets bumped to 45186 on every data row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}

# Rows that carry species/hyperlink data (S, T, V, W, X, Y) get a friendly
# display label (the report id from column A) added as the HYPERLINK()
# function's second argument.
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($label)) {
        continue
    }

    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $r)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        $alreadyTagged = $formula.Contains(', "' + $label + '"') -or $formula.Contains(',"' + $label + '"')
        if ($formula.EndsWith(")") -and -not $alreadyTagged) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
